# Add RemoteScanError column and dual logging functionality
# - Add RemoteScanError column for remote UNC access errors
# - Separate local (ScanError) from remote (RemoteScanError) error handling
# - Clear error columns for successful scans (local and remote)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# New header for the remote-scan error column (col I)
$ws.Range("I1").Value2 = "RemoteScanError"

# Row 3 (NOTPRIME / Windows_2019 / commons-collections4 jar): the previous
# local ScanError ("UNC access denied...") no longer applies to this row -
# clear it out.
$ws.Range("H3").ClearContents()

# Row 4 (LPRIME / poi jar): platform label corrected from the underscored
# form to the human readable "Windows Server 2019".
$ws.Range("B4").Value2 = "Windows Server 2019"

# Row 5 (NOTPRIME2 / commons-collections4 jar): same platform label fix,
# and the UNC-access error now belongs to the new RemoteScanError column
# instead of the local ScanError column.
$ws.Range("B5").Value2 = "Windows Server 2019"
$ws.Range("I5").Value2 = "UNC access denied - cannot determine file existence"
$ws.Range("H5").ClearContents()

# Selection follows the newly widened data range.
$ws.Range("E2:I5").Select()
